# Update Mappings 22 Ontologies
# - fix metadata4Ing -> metadata4ing casing in the header row
# - add a new column F "metadata4ing_DEF" with a definition string per row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header casing: metadata4Ing -> metadata4ing
$ws.Range("D1").Value = "metadata4ing_IRI"
$ws.Range("E1").Value = "metadata4ing_DESC"

# New column F: copy header formatting (bold/centered/bordered) from B1, then set text
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "metadata4ing_DEF"

$ws.Range("F2").Value = '[''p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]'', locstr("Process, i.e., a physical entity with a temporal evolution that ''has a meaning for the ontologist''", ''en'')]'
$ws.Range("F3").Value = '[''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]'', ''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])'']'
$ws.Range("F4").Value = '[locstr(''A role is the function of an entity or agent with respect to an activity, in the context of a usage, generation, invalidation, association, start, and end.'', ''en'')]'
